$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D (shifts old D:K to E:L) for the new FY2018 (2018-12-31) data column.
$ws.Columns("D:D").Insert()

# The freshly inserted column has no number formatting; copy formats from the
# (now shifted) original column E so dates/numbers in the new column D render
# the same way as the rest of the financial-year columns.
$ws.Columns("E:E").Copy()
$ws.Columns("D:D").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D with the latest fiscal-year figures.
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 57800
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 14700
$ws.Range("D18").Value = 43100
$ws.Range("D20").Value = -30700
$ws.Range("D21").Value = 14200
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 12400
$ws.Range("D24").Value = 2200
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 10200
$ws.Range("D27").Value = 10200
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 500
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 30700
$ws.Range("D33").Value = 10700
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 10700
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 23700
$ws.Range("D42").Value = 30600
$ws.Range("D43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 34900
$ws.Range("D49").Value = 2900
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 0
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 1481700
$ws.Range("D57").Value = 1000
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 0
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 49500
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 1382300
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 54100
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 99400
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 10700
$ws.Range("D83").Value = 1800
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 16300
$ws.Range("D91").Value = -2300
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -78600
$ws.Range("D96").Value = -2000
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 41500
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -20800

$ws.Columns("A:L").AutoFit()
